$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the Status column data range (H2:H11) as the last user selection,
# then clear its contents (the "Status" header in H1 remains untouched).
$range = $ws.Range("H2:H11")
$range.Select()
$range.ClearContents()
